# Update "想去人数" (F) and "最低票价" (G) figures across sheets,
# reflecting the data refresh recorded in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "展览" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value  = 366    # F2:  365 -> 366
$ws1.Cells.Item(3, 7).Value  = 50     # G3:   40 -> 50
$ws1.Cells.Item(4, 6).Value  = 10818  # F4:  10799 -> 10818
$ws1.Cells.Item(4, 7).Value  = 75     # G4:   65 -> 75
$ws1.Cells.Item(6, 6).Value  = 978    # F6:  976 -> 978
$ws1.Cells.Item(7, 6).Value  = 165    # F7:  164 -> 165
$ws1.Cells.Item(9, 6).Value  = 8293   # F9:  8289 -> 8293
$ws1.Cells.Item(11, 6).Value = 468    # F11: 467 -> 468
$ws1.Cells.Item(12, 6).Value = 507    # F12: 422 -> 507
$ws1.Cells.Item(13, 6).Value = 218    # F13: 217 -> 218
$ws1.Cells.Item(18, 6).Value = 29     # F18: 26 -> 29
$ws1.Cells.Item(19, 6).Value = 784    # F19: 782 -> 784
$ws1.Cells.Item(20, 6).Value = 131    # F20: 130 -> 131
$ws1.Cells.Item(23, 6).Value = 111    # F23: 109 -> 111
$ws1.Cells.Item(24, 6).Value = 1774   # F24: 1772 -> 1774

# --- Sheet 2: "演出" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value = 31      # F2: 30 -> 31

# --- Sheet 4: "全部类型" (mirrors sheet1 rows 2-24, plus sheet2 row as row 25) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value  = 366    # F2:  365 -> 366
$ws4.Cells.Item(3, 7).Value  = 50     # G3:   40 -> 50
$ws4.Cells.Item(4, 6).Value  = 10818  # F4:  10799 -> 10818
$ws4.Cells.Item(4, 7).Value  = 75     # G4:   65 -> 75
$ws4.Cells.Item(6, 6).Value  = 978    # F6:  976 -> 978
$ws4.Cells.Item(7, 6).Value  = 165    # F7:  164 -> 165
$ws4.Cells.Item(9, 6).Value  = 8293   # F9:  8289 -> 8293
$ws4.Cells.Item(11, 6).Value = 468    # F11: 467 -> 468
$ws4.Cells.Item(12, 6).Value = 507    # F12: 422 -> 507
$ws4.Cells.Item(13, 6).Value = 218    # F13: 217 -> 218
$ws4.Cells.Item(18, 6).Value = 29     # F18: 26 -> 29
$ws4.Cells.Item(19, 6).Value = 784    # F19: 782 -> 784
$ws4.Cells.Item(20, 6).Value = 131    # F20: 130 -> 131
$ws4.Cells.Item(23, 6).Value = 111    # F23: 109 -> 111
$ws4.Cells.Item(24, 6).Value = 1774   # F24: 1772 -> 1774
$ws4.Cells.Item(25, 6).Value = 31     # F25: 30 -> 31
